$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Ativação:" date from 01/01/2012 to 01/01/2023 (rows 8 and 13 share the value) ---
# A direct .Value assignment of a date-shaped string gets auto-coerced into a real
# Excel date serial (and mints a brand new number-format style). To keep this a
# literal text value (and keep the original style s="2"/s="3"), stage the text in a
# scratch cell via a text formula, copy it, and paste-special VALUES ONLY into the
# target cells - that carries over the literal text without re-triggering date
# detection and without touching the destination cell's existing style.
$ws.Range("Z1").Formula = '="01/01/2023"'
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# --- Row 11 ("Objectives:") gains the English objective text in B/C ---
$ws.Range("B11").Value = "To present the concepts of spintronics and the potential applications in quantum computing."
$ws.Range("C11").Value = "To present the concepts of spintronics and the potential applications in quantum computing."
# New cells default to an inherited style that doesn't match the sheet's B/C column
# style, so copy the formatting from a neighboring row that already has it right.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- Row 14 ("Short syllabus:") gains the short syllabus text in B/C ---
$ws.Range("B14").Value = "Introduction to nanotechnology. Metal spintronics. Semiconductor spintronics. Spintronics devices. Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms. decoherence. Quantum Dots. Kane transistor. Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("C14").Value = "Introduction to nanotechnology. Metal spintronics. Semiconductor spintronics. Spintronics devices. Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms. decoherence. Quantum Dots. Kane transistor. Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# --- Row 16 ("Syllabus:") gains the full syllabus text in B/C ---
$ws.Range("B16").Value = "Introduction to nanotechnology.Metal spintronics. Semiconductor SpintronicsSpintronics devices.Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms.decoherence. Quantum Dots.Kane transistor.Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("C16").Value = "Introduction to nanotechnology.Metal spintronics. Semiconductor SpintronicsSpintronics devices.Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms.decoherence. Quantum Dots.Kane transistor.Introduction to models of consciousness: is the brain a quantum computer?"
$ws.Range("B18").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C18").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$excel.CutCopyMode = 0
